$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from an existing header cell (A1) to the new header cells
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Set header values for the new columns
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Set boolean outlier flag values for rows 2-8
$values = @(
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $true,  $true),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false),
    @($false, $false, $false)
)

for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i][0]
    $ws.Cells.Item($row, 7).Value = $values[$i][1]
    $ws.Cells.Item($row, 8).Value = $values[$i][2]
}
